# cryptoVirtualBase.xlsx — "Add files via upload" re-run:
#  - bump the "dollars added per buy" counter in D2 from 1 to 2
#  - move the live selection to D2
#  - log 8 new price-check rows (19-27), rows 19-23 keep the
#    "recent / highlighted" centered look, rows 24-27 are the next
#    purchase batch using the sheet's plain/general formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2: combien de dollars j'ajoute a chaque achat: 1 -> 2 ---
$ws.Range("D2").Value = 2

# --- selection moves to D2 ---
$ws.Range("D2").Select()

function Set-PriceRow {
    param($RowNum, $Price, $Qty, $DateText, $TimeText, $Centered)
    $rng = $ws.Range("A$($RowNum):D$($RowNum)")
    $ws.Range("A$RowNum").Value = $Price
    $ws.Range("B$RowNum").Value = $Qty
    $ws.Range("C$RowNum").Value = $DateText
    $ws.Range("D$RowNum").Value = $TimeText
    if ($Centered) {
        $rng.HorizontalAlignment = -4108   # xlCenter - matches style used by rows 3-19
        $rng.EntireRow.RowHeight = 12.75   # pin the explicit row height like rows 1-19
    } else {
        $rng.HorizontalAlignment = 1       # xlGeneral - matches the plain trailing rows
    }
}

# Rows 19 & 20 already existed; they pick up the centered styling too
Set-PriceRow "19" 1.0587 1 "18/11/2025" "01:01:37" $true
Set-PriceRow "20" 1.0534 1 "19/11/2025" "19:01:50" $true

# New rows 21-23: still the centered/highlighted style
Set-PriceRow "21" 0.9909             1 "20/11/2025" "19:01:26" $true
Set-PriceRow "22" 0.9906             1 "20/11/2025" "20:14:58" $true
Set-PriceRow "23" 0.9762999999999999 1 "21/11/2025" "07:01:07" $true

# New rows 24-27: next purchase batch (qty 2), plain/general formatting
Set-PriceRow "24" 0.9169 2 "21/11/2025" "14:01:37" $false
Set-PriceRow "25" 0.9135 2 "22/11/2025" "01:00:38" $false
Set-PriceRow "26" 0.8812 2 "22/11/2025" "07:00:48" $false
Set-PriceRow "27" 0.8778 2 "22/11/2025" "14:01:08" $false
